$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the interval count input cell; dependent formulas recalculate automatically.
$ws.Range("C4").Value = 20

# Move the active selection as recorded in the sheet view.
$ws.Range("B27").Select()
